$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A85").Value = 90124
